$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D holds numeric-looking text (e.g. "8.07"); Excel would otherwise
# auto-convert these to Number cells. Force Text by switching the number format
# to Text ("@") before writing the value, then restore the original style so the
# cells style index is unchanged (only the value differs from the original).
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") '67.735.77'
$ws.Range("E2").Value = '  -1.30%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.429.05'
$ws.Range("E3").Value = '  -0.92%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
Set-TextValue $ws.Range("D5") '552.61'
$ws.Range("E5").Value = '  -0.80%  '

# Row 6
Set-TextValue $ws.Range("D6") '159.61'
$ws.Range("E6").Value = '  -0.88%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("E8").Value = '  +1.27%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.160'
$ws.Range("E9").Value = '  +7.45%  '

# Row 10
$ws.Range("E10").Value = '  -0.60%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.329'
$ws.Range("E11").Value = '  -1.25%  '

# Row 12
Set-TextValue $ws.Range("D12") '4.77'
$ws.Range("E12").Value = '  -0.16%  '

# Row 13
Set-TextValue $ws.Range("D13") '67.681.21'
$ws.Range("E13").Value = '  -1.29%  '

# Row 14
Set-TextValue $ws.Range("D14") '0.0000170'
$ws.Range("E14").Value = '  +1.73%  '

# Row 15
Set-TextValue $ws.Range("D15") '22.95'
$ws.Range("E15").Value = '  -1.50%  '

# Row 16
$ws.Range("E16").Value = '  -2.96%  '

# Row 17
Set-TextValue $ws.Range("D17") '333.44'
$ws.Range("E17").Value = '  -2.08%  '

# Row 18
Set-TextValue $ws.Range("D18") '6.85'
$ws.Range("E18").Value = '  -1.68%  '

# Row 20
$ws.Range("E20").Value = '  +0.10%  '

# Row 21
$ws.Range("E21").Value = '  +0.53%  '

# Row 22
Set-TextValue $ws.Range("D22") '66.21'
$ws.Range("E22").Value = '  -0.50%  '

# Row 23
Set-TextValue $ws.Range("D23") '3.62'
$ws.Range("E23").Value = '  -0.54%  '

# Row 24
Set-TextValue $ws.Range("D24") '8.07'
$ws.Range("E24").Value = '  -0.07%  '

# Row 25
Set-TextValue $ws.Range("D25") '0.0₃0803'
$ws.Range("E25").Value = '  -1.07%  '

# Row 26
Set-TextValue $ws.Range("D26") '7.03'
$ws.Range("E26").Value = '  -1.43%  '

# Row 27
$ws.Range("E27").Value = '  -0.02%  '

# Row 28
$ws.Range("B28").Value = 'Bittensor'
$ws.Range("C28").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D28") '416.26'
$ws.Range("E28").Value = '  -4.24%  '

# Row 29
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D29") '1.13'
$ws.Range("E29").Value = '  -0.12%  '

# Row 30
Set-TextValue $ws.Range("D30") '1.59'
$ws.Range("E30").Value = '  -0.74%  '

# Row 31
Set-TextValue $ws.Range("D31") '159.06'
$ws.Range("E31").Value = '  +0.95%  '

# Row 32
Set-TextValue $ws.Range("D32") '18.93'

# Row 33
$ws.Range("E33").Value = '  -0.04%  '

# Row 34
Set-TextValue $ws.Range("D34") '17.85'
$ws.Range("E34").Value = '  +0.37%  '

# Row 35
$ws.Range("E35").Value = '  -2.63%  '

# Row 36
$ws.Range("E36").Value = '  -2.30%  '

# Row 37
$ws.Range("E37").Value = '  -3.11%  '

# Row 38
$ws.Range("E38").Value = '  +0.11%  '

# Row 39
$ws.Range("E39").Value = '  -1.30%  '

# Row 40
$ws.Range("E40").Value = '  -2.22%  '

# Row 41
$ws.Range("E41").Value = '  -0.08%  '

# Row 42
Set-TextValue $ws.Range("D42") '129.64'
$ws.Range("E42").Value = '  -1.57%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.0707'
$ws.Range("E43").Value = '  -0.98%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.477'
$ws.Range("E44").Value = '  -0.48%  '

# Row 45
Set-TextValue $ws.Range("D45") '0.554'
$ws.Range("E45").Value = '  -0.83%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.0913'
$ws.Range("E46").Value = '  +0.78%  '

# Row 47
$ws.Range("E47").Value = '  +0.39%  '

# Row 48
Set-TextValue $ws.Range("D48") '1.33'
$ws.Range("E48").Value = '  -6.71%  '

# Row 49
$ws.Range("E49").Value = '  -1.69%  '

# Row 50
$ws.Range("E50").Value = '  +2.70%  '

# Row 51
$ws.Range("E51").Value = '  +0.27%  '
